$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 29 (shifts existing rows 29..119 down to 30..120)
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row with the data for day 28 of July/2025
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = 24381.97
$ws.Cells.Item(29, 3).Value = 7
$ws.Cells.Item(29, 4).Value = 2025
$ws.Cells.Item(29, 5).Value = "07/2025"
